$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / add cell values ---
$ws.Range("A1").Value = "Take over"
$ws.Range("B1").Value = "Addition"
$ws.Range("A3").Value = "Pats, Hugs and Pokes"
$ws.Range("B3").Value = "Add multiple Users as parameter"
$ws.Range("C3").Value = "high5"
$ws.Range("D3").Value = " not more then 1 mention of same person"
$ws.Range("A4").Value = "Dynamic prefix"
$ws.Range("A5").Value = "Self-Assignable Roles"
$ws.Range("A6").Value = "Reminders"
$ws.Range("B6").Value = "actually make it work now :^)"
$ws.Range("A7").Value = "Admin / Mod"
$ws.Range("B7").Value = "after couple warnings auto actions dont work"
$ws.Range("C7").Value = "add reason to official discord log"
$ws.Range("D7").Value = "purge / clear commands"
$ws.Range("A8").Value = "Searches"
$ws.Range("B8").Value = "image search?"
$ws.Range("A9").Value = "League Of Legends"
$ws.Range("A10").Value = "Music"
$ws.Range("B10").Value = "repeat"
$ws.Range("C10").Value = "save playlist maybe? VIP ONLY?"
$ws.Range("A11").Value = "AFK / Away"
$ws.Range("B11").Value = "Dont invoke if the message starts with the Guild's prefix"
$ws.Range("A12").Value = "Custom Member Join / Leave Announcements"
$ws.Range("A13").Value = "Tags"
$ws.Range("B13").Value = "pagify taglist!"
$ws.Range("A14").Value = "Marriages"
$ws.Range("B14").Value = "support divorces with ID only => even without connection to sora"
$ws.Range("C14").Value = "marriage limit!!! "
$ws.Range("A15").Value = "Blacklist"
$ws.Range("A16").Value = "Starboard"
$ws.Range("B16").Value = "Minimum star count per guild setting"
$ws.Range("A17").Value = "Profile and EP"
$ws.Range("A18").Value = "Help"
$ws.Range("A19").Value = "Info"
$ws.Range("B19").Value = "guild doesnt require pic"
$ws.Range("C19").Value = "guild doesnt require emojis"
$ws.Range("A20").Value = "Miscellaneous and Fun"
$ws.Range("A21").Value = "Changelog"
$ws.Range("A23").Value = "http://git.argus.moe/serenity/SoraBot/wikis/sora-help"

# --- Clear the cell whose old content moved away (Changelog moved from A22 to A21) ---
$ws.Range("A22").ClearContents()

# --- Re-apply the "Gut" (good/green) highlight style to the relocated header rows ---
$ws.Range("A3:D3").Style = "Gut"
$ws.Range("A4").Style = "Gut"
$ws.Range("A11:B11").Style = "Gut"

# --- Column C got a bit wider ---
# (The host's ColumnWidth setter quantizes to the nearest 1/6 character width,
#  so 36.5 is the closest input that lands on the target stored width of ~37.29.)
$ws.Range("C1").ColumnWidth = 36.5

# --- Add the hyperlink on the Changelog wiki URL cell (also applies the built-in Link style) ---
$ws.Hyperlinks.Add($ws.Range("A23"), "http://git.argus.moe/serenity/SoraBot/wikis/sora-help") | Out-Null

# --- Restore the selected cell shown when the workbook is opened ---
$ws.Range("I3").Select()
